$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect to apply the data refresh,
# then restore protection so the sheet remains locked as before.
$ws.Unprotect()

# Update confidentiality disclaimer date from 2021-04-29 to 2021-04-30
$ws.Range("A80").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-30 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-77 with refreshed model data
$ws.Range("D2").Value2 = 0.06615001172046582
$ws.Range("E2").Value2 = -0.01513335331135734
$ws.Range("D3").Value2 = 0.04095980085929946
$ws.Range("E3").Value2 = -0.001120614407817189
$ws.Range("D4").Value2 = 0.03277447201915092
$ws.Range("E4").Value2 = -0.001306878935487643
$ws.Range("D5").Value2 = 0.0307431842004518
$ws.Range("E5").Value2 = -0.01393228067212693
$ws.Range("D6").Value2 = 0.02823342573959035
$ws.Range("E6").Value2 = -0.01640783028803572
$ws.Range("D7").Value2 = 0.0244155722012907
$ws.Range("E7").Value2 = -0.008892325536439372
$ws.Range("D8").Value2 = 0.1811521709937732
$ws.Range("E8").Value2 = -0.02571785268414484
$ws.Range("D9").Value2 = 0.02389560657125484
$ws.Range("E9").Value2 = -0.008952496954933054
$ws.Range("D10").Value2 = 0.02188976638966868
$ws.Range("E10").Value2 = 0.006867406233491691
$ws.Range("D11").Value2 = 0.02240477622020496
$ws.Range("E11").Value2 = -0.007043094186462762
$ws.Range("D12").Value2 = 0.02041018491684802
$ws.Range("E12").Value2 = 0.003723088544757891
$ws.Range("D13").Value2 = 0.0197265205763576
$ws.Range("E13").Value2 = -0.01410848941863296
$ws.Range("D14").Value2 = 0.01679929500526916
$ws.Range("E14").Value2 = -0.01068791294209104
$ws.Range("D15").Value2 = 0.01721455953794077
$ws.Range("E15").Value2 = -0.008609108107367147
$ws.Range("D16").Value2 = 0.01574780021303173
$ws.Range("E16").Value2 = -0.02893222506393867
$ws.Range("D17").Value2 = 0.01427891697405144
$ws.Range("E17").Value2 = 0.002005299720690301
$ws.Range("D18").Value2 = 0.01415352738036352
$ws.Range("E18").Value2 = 0.001528414617201657
$ws.Range("D19").Value2 = 0.01425622255888271
$ws.Range("E19").Value2 = -0.01344420503171373
$ws.Range("D20").Value2 = 0.01298199211135836
$ws.Range("E20").Value2 = -0.02884289107567017
$ws.Range("D21").Value2 = 0.01220330228205049
$ws.Range("E21").Value2 = 0.002233567326100916
$ws.Range("D22").Value2 = 0.01264834094346167
$ws.Range("E22").Value2 = 0.005500946884299696
$ws.Range("D23").Value2 = 0.01223846485945231
$ws.Range("E23").Value2 = -0.01771435917213005
$ws.Range("D24").Value2 = 0.01263492881941914
$ws.Range("E24").Value2 = -0.006848462208940309
$ws.Range("D25").Value2 = 0.01151657006567937
$ws.Range("E25").Value2 = -0.005133091535009404
$ws.Range("D26").Value2 = 0.009697948976299773
$ws.Range("E26").Value2 = -0.02259822848058157
$ws.Range("D27").Value2 = 0.01008418668333095
$ws.Range("E27").Value2 = -0.03201397881336099
$ws.Range("D28").Value2 = 0.01029980329330497
$ws.Range("E28").Value2 = -0.02081949058693244
$ws.Range("D29").Value2 = 0.01035789627632789
$ws.Range("E29").Value2 = -0.001742947517913707
$ws.Range("D30").Value2 = 0.01013736319859634
$ws.Range("E30").Value2 = -0.006130208737487264
$ws.Range("D31").Value2 = 0.008939790316347326
$ws.Range("E31").Value2 = -0.01286890871654078
$ws.Range("D32").Value2 = 0.01020422716009869
$ws.Range("E32").Value2 = -0.004432624113475225
$ws.Range("D33").Value2 = 0.009243431365052082
$ws.Range("E33").Value2 = 0.008199581297976266
$ws.Range("D34").Value2 = 0.008963389361583448
$ws.Range("E34").Value2 = -0.005160339107998624
$ws.Range("D35").Value2 = 0.009076586115232708
$ws.Range("E35").Value2 = -0.001343328855570647
$ws.Range("D36").Value2 = 0.008567161392068978
$ws.Range("E36").Value2 = -0.007685315263201487
$ws.Range("D37").Value2 = 0.008555833850355641
$ws.Range("E37").Value2 = -0.01441180526823882
$ws.Range("D38").Value2 = 0.00798827681242693
$ws.Range("E38").Value2 = 0.04791728212703106
$ws.Range("D39").Value2 = 0.008501949363733163
$ws.Range("E39").Value2 = 0.001295336787564771
$ws.Range("D40").Value2 = 0.007988670129847534
$ws.Range("E40").Value2 = -0.03582787652011232
$ws.Range("D41").Value2 = 0.007335999202100549
$ws.Range("E41").Value2 = -0.02163889425035925
$ws.Range("D42").Value2 = 0.007856594140009377
$ws.Range("E42").Value2 = -0.0271937202130641
$ws.Range("D43").Value2 = 0.007758658102279475
$ws.Range("E43").Value2 = -0.001054435218136196
$ws.Range("D44").Value2 = 0.007337887125719438
$ws.Range("E44").Value2 = -0.004802641452798961
$ws.Range("D45").Value2 = 0.007306107078134796
$ws.Range("E45").Value2 = 0.01257563685695229
$ws.Range("D46").Value2 = 0.007685737052499861
$ws.Range("E46").Value2 = -0.004175878162613555
$ws.Range("D47").Value2 = 0.006955110611989558
$ws.Range("E47").Value2 = 0.01112920738327894
$ws.Range("D48").Value2 = 0.007146577532338618
$ws.Range("E48").Value2 = -0.006075949367088662
$ws.Range("D49").Value2 = 0.006762896388541352
$ws.Range("E49").Value2 = -0.001221320771176759
$ws.Range("D50").Value2 = 0.0066082833105027
$ws.Range("E50").Value2 = 0.00039282440749
$ws.Range("D51").Value2 = 0.00626275395650383
$ws.Range("E51").Value2 = 0.002813557831802038
$ws.Range("D52").Value2 = 0.006507987368249187
$ws.Range("E52").Value2 = -0.009693951554416769
$ws.Range("D53").Value2 = 0.005472618590256441
$ws.Range("E53").Value2 = 0.01164294954721856
$ws.Range("D54").Value2 = 0.005732994722694975
$ws.Range("E54").Value2 = 0.00312843029637766
$ws.Range("D55").Value2 = 0.005986487800272974
$ws.Range("E55").Value2 = -0.03133931211195418
$ws.Range("D56").Value2 = 0.005650059811386835
$ws.Range("E56").Value2 = -0.007910817437212558
$ws.Range("D57").Value2 = 0.006410129994003404
$ws.Range("E57").Value2 = 0.0006872177498529197
$ws.Range("D58").Value2 = 0.005525952432490075
$ws.Range("E58").Value2 = -0.004783054321831237
$ws.Range("D59").Value2 = 0.005227660500705507
$ws.Range("E59").Value2 = -0.03521126760563398
$ws.Range("D60").Value2 = 0.004827420693500895
$ws.Range("E60").Value2 = 0.009125277017338007
$ws.Range("D61").Value2 = 0.004650545849456168
$ws.Range("E61").Value2 = -0.03488696622941667
$ws.Range("D62").Value2 = 0.004586081124219499
$ws.Range("E62").Value2 = -0.01037735849056598
$ws.Range("D63").Value2 = 0.004173255159555623
$ws.Range("E63").Value2 = 0.005202442886224734
$ws.Range("D64").Value2 = 0.003957717213065719
$ws.Range("E64").Value2 = -0.007473366194943543
$ws.Range("D65").Value2 = 0.003955514635510348
$ws.Range("E65").Value2 = -0.01038103571712667
$ws.Range("D66").Value2 = 0.003648412393504294
$ws.Range("E66").Value2 = 0.001940491591203353
$ws.Range("D67").Value2 = 0.003836496784036178
$ws.Range("E67").Value2 = -0.01070308174940016
$ws.Range("D68").Value2 = 0.00349336666630298
$ws.Range("E68").Value2 = -0.01942173883672227
$ws.Range("D69").Value2 = 0.003521292203165723
$ws.Range("E69").Value2 = -0.005361451166115527
$ws.Range("D70").Value2 = 0.003189489627145862
$ws.Range("E70").Value2 = -0.01465002712967989
$ws.Range("D71").Value2 = 0.003160502133247495
$ws.Range("E71").Value2 = -0.02075788687698332
$ws.Range("D72").Value2 = 0.002410996456548294
$ws.Range("E72").Value2 = -0.02057129806358993
$ws.Range("D73").Value2 = 0.002001985670864258
$ws.Range("E73").Value2 = 0.008781925343811503
$ws.Range("D74").Value2 = 0.002029871875984941
$ws.Range("E74").Value2 = -0.01501676064252366
$ws.Range("D75").Value2 = 0.001498303382041317
$ws.Range("E75").Value2 = 0.02756339581036382
$ws.Range("D76").Value2 = 0.001492324957248167
$ws.Range("E76").Value2 = -0.03732012018343789
$ws.Range("D77").Value2 = 1
$ws.Range("E77").Value2 = -0.01116157749454061

# Restore sheet protection
$ws.Protect()
